# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker data table (rows 16-33, columns B:G) is reordered so that all
# "1801" period rows come first (rows 16-24) followed by all "1802" period
# rows (rows 25-33), keeping the same worker order within each period block.
# Also corrects FABIAN PITALUA ZARZA's Salario Basico (col G) from 1142000
# to 1800000 for both periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired state for B16:G33 -> (TipoDoc, NumDoc, Nombre, Periodo, ValorMora, SalarioBasico)
$data = @(
    @("CC", "73182666",   "FABIAN PITALUA ZARZA",              "1801", 72000,  1800000),
    @("CC", "1143343026", "ANDRES FELIPE VASQUEZ MEJIA",       "1801", 96000,  2400000),
    @("CC", "33069585",   "MARIA CLARA URIBE AGUILAR",         "1801", 73771,  1844292),
    @("CC", "1143374517", "DARWIN CARIAGA GARCIA",             "1801", 36000,  900000),
    @("CC", "1143349287", "YENIFER PAOLA TATAR RODRIGUEZ",     "1801", 40000,  1000000),
    @("CE", "362441",     "LUIS FRANCISCO SAGARZAZU RODRIGUEZ","1801", 29509,  737717),
    @("CE", "501276",     "MARCOS JOSE BORGES RAMOS",          "1801", 96000,  2400000),
    @("CC", "16787235",   "HECTOR FABIO FIGUEROA SOJET",       "1801", 240000, 6000000),
    @("CC", "1127618941", "LEONARDO JOSE ROJAS LARA",          "1801", 96000,  2400000),
    @("CC", "73182666",   "FABIAN PITALUA ZARZA",              "1802", 72000,  1800000),
    @("CC", "1143343026", "ANDRES FELIPE VASQUEZ MEJIA",       "1802", 96000,  2400000),
    @("CC", "33069585",   "MARIA CLARA URIBE AGUILAR",         "1802", 73771,  1844292),
    @("CC", "1143374517", "DARWIN CARIAGA GARCIA",             "1802", 36000,  900000),
    @("CC", "1143349287", "YENIFER PAOLA TATAR RODRIGUEZ",     "1802", 14667,  1000000),
    @("CE", "362441",     "LUIS FRANCISCO SAGARZAZU RODRIGUEZ","1802", 29509,  737717),
    @("CE", "501276",     "MARCOS JOSE BORGES RAMOS",          "1802", 96000,  2400000),
    @("CC", "16787235",   "HECTOR FABIO FIGUEROA SOJET",       "1802", 240000, 6000000),
    @("CC", "1127618941", "LEONARDO JOSE ROJAS LARA",          "1802", 96000,  2400000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
